$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new row with a UTF-8 (Hindi) greeting in B4, matching the style of
# the surrounding data cells (centered horizontal alignment == style index 1).
$ws.Range("B4").Value = "नमस्ते"
$ws.Range("B4").HorizontalAlignment = -4108

# Select B4 so the saved sheetView records it as the active cell.
$ws.Range("B4").Select()

# Add a new defined name with a UTF-8 (Greek) name pointing at Sheet1!$B$4.
# Names.Add() chokes when the `name` argument itself contains non-ASCII
# characters, so create it with a placeholder ASCII name first and rename
# it afterwards.
$n = $wb.Names.Add("TempGreekName", '=Sheet1!$B$4')
$n.Name = 'Χαιρετισμός'
